# Weekly update: insert two new observation rows for the latest week
# (2023-04-25) right after the current "latest" row (row 1193, dated
# 2023-03-27), pushing the rest of the historical records down by two
# rows, and append the two rows that used to be the tail of the table
# (the insert naturally carries those down to the new end of range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 1194:1195 -- this shifts the existing rows
# 1194..1263 down to 1196..1265 and copies formatting from the row above,
# matching Excel's native "Insert Copied Cells"/"Insert Rows" behaviour.
$ws.Range("A1194:R1195").Insert()

# --- New row 1194 ---
$ws.Range("A1194").Value = 6
$ws.Range("B1194").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1194").Value = "Metropolitana"
$ws.Range("D1194").Value = 45041
$ws.Range("E1194").Value = 13
$ws.Range("F1194").Value = 100112040
$ws.Range("G1194").Value = "Cilantro"
$ws.Range("H1194").Value = "Sin especificar"
$ws.Range("I1194").Value = "Primera"
$ws.Range("J1194").Value = 770
$ws.Range("K1194").Value = 5000
$ws.Range("L1194").Value = 5500
$ws.Range("M1194").Value = 5227
$ws.Range("N1194").Value = "`$/caja 36 atados"
$ws.Range("O1194").Value = "Región Metropolitana"
$ws.Range("P1194").Value = 145
$ws.Range("Q1194").Value = 36
$ws.Range("R1194").Value = "Hortaliza"

# --- New row 1195 ---
$ws.Range("A1195").Value = 6
$ws.Range("B1195").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1195").Value = "Metropolitana"
$ws.Range("D1195").Value = 45041
$ws.Range("E1195").Value = 13
$ws.Range("F1195").Value = 100112040
$ws.Range("G1195").Value = "Cilantro"
$ws.Range("H1195").Value = "Sin especificar"
$ws.Range("I1195").Value = "Primera"
$ws.Range("J1195").Value = 450
$ws.Range("K1195").Value = 8000
$ws.Range("L1195").Value = 9000
$ws.Range("M1195").Value = 8400
$ws.Range("N1195").Value = "`$/docena de atados"
$ws.Range("O1195").Value = "Región Metropolitana"
$ws.Range("P1195").Value = 2800
$ws.Range("Q1195").Value = 3
$ws.Range("R1195").Value = "Hortaliza"
